# Atualizações dados 17/07 19h
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDt = 45490.81285879629

# Update the dt_insertion (column H) timestamp for every data row (2-21)
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 8).Value2 = $newDt
}

# Row 12: points (F) 20 -> 21, matches (G) 16 -> 17
$ws.Cells.Item(12, 6).Value2 = 21
$ws.Cells.Item(12, 7).Value2 = 17

# Row 19: position (E) 18 -> 19
$ws.Cells.Item(19, 5).Value2 = 19

# Row 20: position (E) 19 -> 18, points (F) 11 -> 12, matches (G) 16 -> 17
$ws.Cells.Item(20, 5).Value2 = 18
$ws.Cells.Item(20, 6).Value2 = 12
$ws.Cells.Item(20, 7).Value2 = 17
